$d = $word.ActiveDocument

# Color used for highlighted metrics: RGB(0x2C, 0x3E, 0x50) stored as
# Word's BGR-packed Long (r + g*256 + b*65536).
$metricColor = 5258796

function Apply-MetricHighlight($ParagraphText, $Metrics) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.TrimEnd() -eq $ParagraphText) {
            foreach ($metric in $Metrics) {
                $fr = $p.Range.Duplicate
                $fr.Find.Execute($metric, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
                $fr.Font.Bold = 1
                $fr.Font.Color = $metricColor
            }
        }
    }
}

Apply-MetricHighlight "• Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%" @("23%", "64%")

Apply-MetricHighlight "• Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from ±4.2% to ±2.1%" @("87%", "71%", "±4.2%", "±2.1%")

Apply-MetricHighlight "• Wrote RFP and analyzed bids from 1,200 vendors for research platform development" @("1,200")

Apply-MetricHighlight "• Created comprehensive meta-analysis framework handling millions of survey responses that became the `$400M Polling Consortium Database at The Analyst Institute, now valued at `$1B+" @("`$400M", "`$1B")

Apply-MetricHighlight "• Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M" @("73.5%", "`$4.7M")

Apply-MetricHighlight "• Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%" @("87%", "71%")
